$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 363.68
$ws.Range("I92").Value = 327.5909
$ws.Range("K92").Value = 327.5909
$ws.Range("M92").Value = 920.4091000000001
$ws.Range("H111").Value = 52631930
$ws.Range("I111").Value = 83333690
$ws.Range("J111").Value = 341.42856
$ws.Range("K111").Value = 250001070
$ws.Range("L111").Value = 1024.28568
$ws.Range("M111").Value = -249998003
$ws.Range("N111").Value = -7158.28568
$ws.Range("H132").Value = 2760.3396
$ws.Range("I132").Value = 2555.081
$ws.Range("K132").Value = 7665.243
$ws.Range("M132").Value = -5135.243
$ws.Range("H135").Value = 451
$ws.Range("I135").Value = 492.78946
$ws.Range("J135").Value = 252.5
$ws.Range("K135").Value = 4435.105140000001
$ws.Range("L135").Value = 2272.5
$ws.Range("M135").Value = -1900.105140000001
$ws.Range("N135").Value = -7342.5
$ws.Range("H138").Value = 2073.5925
$ws.Range("I138").Value = 975.8095
$ws.Range("J138").Value = 5915.8335
$ws.Range("K138").Value = 2927.4285
$ws.Range("L138").Value = 17747.5005
$ws.Range("M138").Value = 2212.5715
$ws.Range("N138").Value = -28027.5005

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1831.25
$ws.Range("I2").Value = 1575.4667
$ws.Range("J2").Value = 2598.6
$ws.Range("K2").Value = 1575.4667
$ws.Range("L2").Value = 2598.6
$ws.Range("M2").Value = -1462.4667
$ws.Range("N2").Value = -2824.6
$ws.Range("H32").Value = 16954302
$ws.Range("I32").Value = 18871412
$ws.Range("K32").Value = 18871412
$ws.Range("M32").Value = -18871125
$ws.Range("H45").Value = 1842.8
$ws.Range("I45").Value = 1400
$ws.Range("J45").Value = 2507
$ws.Range("K45").Value = 1400
$ws.Range("L45").Value = 2507
$ws.Range("M45").Value = -1023
$ws.Range("N45").Value = -3261
$ws.Range("H61").Value = 1575.1724
$ws.Range("I61").Value = 1010.5
$ws.Range("J61").Value = 4285.6
$ws.Range("K61").Value = 1010.5
$ws.Range("L61").Value = 4285.6
$ws.Range("M61").Value = -798.5
$ws.Range("N61").Value = -4709.6
$ws.Range("H74").Value = 2401.5386
$ws.Range("I74").Value = 2744.6667
$ws.Range("J74").Value = 1933.6364
$ws.Range("K74").Value = 2744.6667
$ws.Range("L74").Value = 1933.6364
$ws.Range("M74").Value = -1870.6667
$ws.Range("N74").Value = -3681.6364
$ws.Range("H77").Value = 2401.5386
$ws.Range("I77").Value = 2744.6667
$ws.Range("J77").Value = 1933.6364
$ws.Range("K77").Value = 13723.3335
$ws.Range("L77").Value = 9668.182000000001
$ws.Range("M77").Value = -9355.333500000001
$ws.Range("N77").Value = -18404.182
$ws.Range("H97").Value = 541.875
$ws.Range("I97").Value = 361.35715
$ws.Range("J97").Value = 1805.5
$ws.Range("K97").Value = 361.35715
$ws.Range("L97").Value = 1805.5
$ws.Range("M97").Value = 134.64285
$ws.Range("N97").Value = -2797.5
$ws.Range("H116").Value = 1831.25
$ws.Range("I116").Value = 1575.4667
$ws.Range("J116").Value = 2598.6
$ws.Range("K116").Value = 1575.4667
$ws.Range("L116").Value = 2598.6
$ws.Range("M116").Value = 718.5333000000001
$ws.Range("N116").Value = -7186.6
$ws.Range("H132").Value = 1738.4348
$ws.Range("I132").Value = 1620.8846
$ws.Range("J132").Value = 1891.25
$ws.Range("K132").Value = 4862.6538
$ws.Range("L132").Value = 5673.75
$ws.Range("M132").Value = -2332.6538
$ws.Range("N132").Value = -10733.75
$ws.Range("H136").Value = 1575.1724
$ws.Range("I136").Value = 1010.5
$ws.Range("J136").Value = 4285.6
$ws.Range("K136").Value = 3031.5
$ws.Range("L136").Value = 12856.8
$ws.Range("M136").Value = -481.5
$ws.Range("N136").Value = -17956.8
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1831.25
$ws.Range("I3").Value = 1575.4667
$ws.Range("J3").Value = 2598.6
$ws.Range("K3").Value = 1575.4667
$ws.Range("L3").Value = 2598.6
$ws.Range("M3").Value = -1461.4667
$ws.Range("N3").Value = -2826.6
$ws.Range("H134").Value = 1436.9714
$ws.Range("I134").Value = 1291.28
$ws.Range("J134").Value = 1801.2
$ws.Range("K134").Value = 3873.84
$ws.Range("L134").Value = 5403.6
$ws.Range("M134").Value = -1338.84
$ws.Range("N134").Value = -10473.6

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3321.3928
$ws.Range("I31").Value = 2115
$ws.Range("J31").Value = 5868.222
$ws.Range("K31").Value = 2115
$ws.Range("L31").Value = 5868.222
$ws.Range("M31").Value = -1820
$ws.Range("N31").Value = -6458.222
$ws.Range("H34").Value = 3321.3928
$ws.Range("I34").Value = 2115
$ws.Range("J34").Value = 5868.222
$ws.Range("K34").Value = 2115
$ws.Range("L34").Value = 5868.222
$ws.Range("M34").Value = -1913
$ws.Range("N34").Value = -6272.222
$ws.Range("H132").Value = 1985.5143
$ws.Range("I132").Value = 1531.375
$ws.Range("K132").Value = 4594.125
$ws.Range("M132").Value = -2064.125
$ws.Range("H134").Value = 2201.4866
$ws.Range("I134").Value = 1363.3572
$ws.Range("J134").Value = 4809
$ws.Range("K134").Value = 4090.0716
$ws.Range("L134").Value = 14427
$ws.Range("M134").Value = -1555.0716
$ws.Range("N134").Value = -19497

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 3104.8
$ws.Range("I94").Value = 1762
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 5286
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = -4610
$ws.Range("N94").Value = -13352
$ws.Range("H137").Value = 3509.1177
$ws.Range("I137").Value = 3506.8462
$ws.Range("J137").Value = 3516.5
$ws.Range("K137").Value = 10520.5386
$ws.Range("L137").Value = 10549.5
$ws.Range("M137").Value = -5420.5386
$ws.Range("N137").Value = -20749.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2834.5557
$ws.Range("I80").Value = 2728.8572
$ws.Range("J80").Value = 3204.5
$ws.Range("K80").Value = 2728.8572
$ws.Range("L80").Value = 3204.5
$ws.Range("M80").Value = -1730.8572
$ws.Range("N80").Value = -5200.5
$ws.Range("H83").Value = 2834.5557
$ws.Range("I83").Value = 2728.8572
$ws.Range("J83").Value = 3204.5
$ws.Range("K83").Value = 13644.286
$ws.Range("L83").Value = 16022.5
$ws.Range("M83").Value = -8652.286
$ws.Range("N83").Value = -26006.5
$ws.Range("H113").Value = 9246
$ws.Range("I113").Value = 892
$ws.Range("J113").Value = 15511.5
$ws.Range("K113").Value = 892
$ws.Range("L113").Value = 15511.5
$ws.Range("M113").Value = 1278
$ws.Range("N113").Value = -19851.5
$ws.Range("H132").Value = 3472.8545
$ws.Range("I132").Value = 3437.9092
$ws.Range("K132").Value = 10313.7276
$ws.Range("M132").Value = -7783.7276

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1066.3334
$ws.Range("I46").Value = 949.5
$ws.Range("J46").Value = 1300
$ws.Range("K46").Value = 949.5
$ws.Range("L46").Value = 1300
$ws.Range("M46").Value = -761.5
$ws.Range("N46").Value = -1676
$ws.Range("H55").Value = 709.4400000000001
$ws.Range("I55").Value = 664.8095
$ws.Range("J55").Value = 943.75
$ws.Range("K55").Value = 664.8095
$ws.Range("L55").Value = 943.75
$ws.Range("M55").Value = -491.8095
$ws.Range("N55").Value = -1289.75
$ws.Range("H61").Value = 1032.7059
$ws.Range("I61").Value = 972.1667
$ws.Range("J61").Value = 1178
$ws.Range("K61").Value = 972.1667
$ws.Range("L61").Value = 1178
$ws.Range("M61").Value = -770.1667
$ws.Range("N61").Value = -1582
$ws.Range("H75").Value = 34157
$ws.Range("I75").Value = 34157
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 34157
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -33221
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 34157
$ws.Range("I78").Value = 34157
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 102471
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -97791
$ws.Range("N78").ClearContents()
$ws.Range("H113").Value = 1032.7059
$ws.Range("I113").Value = 972.1667
$ws.Range("J113").Value = 1178
$ws.Range("K113").Value = 972.1667
$ws.Range("L113").Value = 1178
$ws.Range("M113").Value = 1197.8333
$ws.Range("N113").Value = -5518
$ws.Range("H122").Value = 4054.0908
$ws.Range("I122").Value = 3927.1428
$ws.Range("J122").Value = 4276.25
$ws.Range("K122").Value = 11781.4284
$ws.Range("L122").Value = 12828.75
$ws.Range("M122").Value = -9331.428400000001
$ws.Range("N122").Value = -17728.75
$ws.Range("H132").Value = 1570.95
$ws.Range("I132").Value = 1730.4819
$ws.Range("J132").Value = 792.05884
$ws.Range("K132").Value = 5191.4457
$ws.Range("L132").Value = 2376.17652
$ws.Range("M132").Value = -2661.4457
$ws.Range("N132").Value = -7436.17652

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 639.5
$ws.Range("I81").Value = 569.4
$ws.Range("J81").Value = 990
$ws.Range("K81").Value = 1138.8
$ws.Range("L81").Value = 1980
$ws.Range("M81").Value = -77.79999999999995
$ws.Range("N81").Value = -4102
$ws.Range("H84").Value = 639.5
$ws.Range("I84").Value = 569.4
$ws.Range("J84").Value = 990
$ws.Range("K84").Value = 5694
$ws.Range("L84").Value = 9900
$ws.Range("M84").Value = -390
$ws.Range("N84").Value = -20508
$ws.Range("H107").Value = 290.61905
$ws.Range("I107").Value = 293
$ws.Range("J107").Value = 289.66666
$ws.Range("K107").Value = 879
$ws.Range("L107").Value = 868.9999799999999
$ws.Range("M107").Value = 1041
$ws.Range("N107").Value = -4708.99998
$ws.Range("H132").Value = 2095.0527
$ws.Range("I132").Value = 1334.4
$ws.Range("J132").Value = 2591.1304
$ws.Range("K132").Value = 4003.2
$ws.Range("L132").Value = 7773.3912
$ws.Range("M132").Value = -1473.2
$ws.Range("N132").Value = -12833.3912
$ws.Range("H136").Value = 3115.255
$ws.Range("I136").Value = 2986.4048
$ws.Range("K136").Value = 8959.214399999999
$ws.Range("M136").Value = -6409.214399999999
